# Template_Database.xlsx edit:
#  - Remove the old "Table1" ListObject (convert back to a plain range)
#  - Replace the door-prize participant data with the new/updated list
#  - Adjust column widths, dimension/selection bookkeeping follows automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Unlist the table (keeps the data, drops xl/tables/table1.xml + the
#        worksheet <tableParts> relationship) ---------------------------------
if ($ws.ListObjects.Count -gt 0) {
    $ws.ListObjects.Item(1).Unlist()
}

# --- 2. Clear out the old data area completely (values + formatting) so the
#        old strings/styles are dropped before we type the new table in ------
$ws.Range("A1:E19").ClearContents()
$ws.Range("A1:E19").ClearFormats()

# --- 3. Header row (KPK / Telephone Number / ISSELECTED first; Name and
#        DOORPRICENAME are (re)typed last further down) ------------------------
$ws.Range("B1").Value = "KPK"
$ws.Range("C1").Value = "Telephone Number"
$ws.Range("D1").Value = "ISSELECTED"

# --- 4. Data rows (Name, KPK code, rank) ---------------------------------------
$names = @(
    "Abraham Naibrohu",
    "Erlenni",
    "Hakim Adni",
    "Tissa Shakira",
    "Syfa",
    "Jocelyn",
    "Kayana",
    "Samuel Dofransrael",
    "Manusa Keren",
    "Manusia Buruk",
    "Ikan Laut",
    "Ikat Air Tawar",
    "Kuda Terbang",
    "Kuda Tanduk",
    "Manusa Biasa",
    "Bukan Orang",
    "Hanya Orang",
    "Orang Utan"
)
$codes = @(
    "I00212","I00213","I00214","I00215","I00216","I00217","I00218","I00219",
    "I00220","I00221","I00222","I00223","I00224","I00225","I00226","I00227",
    "I00228","I00229"
)

# column A (names) filled top-to-bottom first ...
for ($i = 0; $i -lt $names.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $names[$i]
}
# ... then column B (KPK codes) top-to-bottom ...
for ($i = 0; $i -lt $codes.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $codes[$i]
}
# ... then column C (rank numbers)
for ($i = 0; $i -lt $names.Count; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = ($i + 1)
}

# --- 5. Header for E1 + A1 written last (matches upload/reorder pass) ---------
$ws.Range("E1").Value = "DOORPRICENAME"
$ws.Range("A1").Value = "Name"

# --- 6. Column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.15625
$ws.Columns.Item(3).ColumnWidth = 14.75
$ws.Columns.Item(4).ColumnWidth = 9.5
$ws.Columns.Item(5).ColumnWidth = 14

# --- 7. Selection (matches the saved cursor position in the target file) ------
[void]$ws.Range("K20").Select()
